# Applies the crypto price-list refresh described in the commit message:
# "Updated cryptos list on Fri Sep 27 02:33:47 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain text such as "598.94", "1.00", "0.0000186" or
# "64.950.13". Excel auto-coerces such strings into numbers when assigned directly,
# which would lose the exact textual formatting (trailing zeros, "." used as a
# thousands separator, etc.). Mark every Price cell that is about to be rewritten as
# Text first so the new value is stored verbatim, just like the original file.
$priceCellRefs = @(
    "D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13",
    "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D23", "D24", "D25",
    "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36",
    "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47",
    "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $priceCellRefs) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.950.13"
$ws.Range("E2").Value = "  +2.74%  "
$ws.Range("D3").Value = "2.619.60"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "598.94"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").Value = "154.13"
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").Value = "0.117"
$ws.Range("E9").Value = "  +7.83%  "
$ws.Range("D10").Value = "0.404"
$ws.Range("E10").Value = "  +5.05%  "
$ws.Range("D11").Value = "5.73"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "0.153"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").Value = "28.99"
$ws.Range("E13").Value = "  +5.73%  "
$ws.Range("D14").Value = "0.0000186"
$ws.Range("E14").Value = "  +22.36%  "
$ws.Range("D15").Value = "3.089.94"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "64.900.37"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").Value = "2.604.97"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "12.46"
$ws.Range("E18").Value = "  +3.44%  "
$ws.Range("E19").Value = "  +4.98%  "
$ws.Range("D20").Value = "358.64"
$ws.Range("E20").Value = "  +4.45%  "
$ws.Range("D21").Value = "7.33"
$ws.Range("E21").Value = "  +7.82%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "68.46"
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("D24").Value = "9.37"
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("D25").Value = "1.64"
$ws.Range("E25").Value = "  -3.68%  "
$ws.Range("D26").Value = "1.65"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  +2.80%  "
$ws.Range("D28").Value = "8.08"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0943"
$ws.Range("E30").Value = "  +13.18%  "
$ws.Range("D31").Value = "527.63"
$ws.Range("E31").Value = "  -5.47%  "
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").Value = "  +5.18%  "
$ws.Range("D33").Value = "1.80"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").Value = "5.50"
$ws.Range("E34").Value = "  +4.16%  "
$ws.Range("D35").Value = "6.30"
$ws.Range("E35").Value = "  +5.87%  "
$ws.Range("D36").Value = "0.424"
$ws.Range("E36").Value = "  +3.63%  "
$ws.Range("D37").Value = "20.27"
$ws.Range("E37").Value = "  +5.18%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "162.45"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "2.01"
$ws.Range("E39").Value = "  +5.67%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "42.12"
$ws.Range("E42").Value = "  +6.89%  "
$ws.Range("D43").Value = "164.14"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").Value = "4.12"
$ws.Range("E44").Value = "  +3.87%  "
$ws.Range("D45").Value = "0.0610"
$ws.Range("E45").Value = "  +5.53%  "
$ws.Range("D46").Value = "23.02"
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("D47").Value = "2.19"
$ws.Range("E47").Value = "  +7.00%  "
$ws.Range("D48").Value = "0.650"
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("D49").Value = "0.0259"
$ws.Range("E49").Value = "  +5.05%  "
$ws.Range("D50").Value = "0.0977"
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("D51").Value = "19.41"
$ws.Range("E51").Value = "  +2.37%  "
